$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 5 new columns at N:R (before the old "RPO"/"RTO" columns) ---
$ws.Range("N1:R1").EntireColumn.Insert() | Out-Null

# After the insert, the old N/O columns (and everything right of them) shifted
# right by 5: old N6 ("Backup/Recovery", style s17) is now at S6, old O6 (blank,
# style s10) is now at T6, old N7 ("RPO") / O7 ("RTO") are now at S7 / T7.

# Move the "Backup/Recovery" header (value + style) back onto N6.
$ws.Range("S6").Cut($ws.Range("N6")) | Out-Null

# O6:Q6 should share N6's style (s17); R6 should share the old O6 style (now at T6, s10).
$ws.Range("N6").Copy() | Out-Null
$ws.Range("O6:Q6").PasteSpecial(-4122) | Out-Null
$ws.Range("T6").Copy() | Out-Null
$ws.Range("R6").PasteSpecial(-4122) | Out-Null

# N7:R7 should share the old N7/O7 header style (now at S7, s25).
$ws.Range("S7").Copy() | Out-Null
$ws.Range("N7:R7").PasteSpecial(-4122) | Out-Null

# Now remove the old RPO/RTO columns (now shifted to S:T).
$ws.Range("S1:T1").EntireColumn.Delete() | Out-Null

# --- Fill in the new header row (row 7) text ---
$ws.Range("N7").Value = "RPO Goal"
$ws.Range("O7").Value = "RTO Goal"
$ws.Range("P7").Value = "Full Backups"
$ws.Range("Q7").Value = "Log Backups"
$ws.Range("R7").Value = "CHECKDB"

# --- Fill in the new data rows (8, 9, 10) ---
$ws.Range("N8").Value = "1 minute"
$ws.Range("O8").Value = "1 hour"
$ws.Range("P8").Value = "Daily 11PM"
$ws.Range("Q8").Value = "Hourly"
$ws.Range("R8").Value = "Sat 10PM"

$ws.Range("N9").Value = "1 day"
$ws.Range("O9").Value = "1 day"
$ws.Range("P9").Value = "Daily 11PM"
$ws.Range("Q9").Value = "Hourly"
$ws.Range("R9").Value = "Sat 11PM"

$ws.Range("N10").Value = "1 hour"
$ws.Range("O10").Value = "1 day"
$ws.Range("P10").Value = "Weekly Sat 9AM"
$ws.Range("Q10").Value = "N/A"
$ws.Range("R10").Value = "Sun 2AM"

# --- Column widths for the new P (bestFit) / Q (bestFit) columns ---
$ws.Range("P1").ColumnWidth = 14.5
$ws.Range("Q1").ColumnWidth = 11

# --- Update the selection to match the saved workbook state ---
$ws.Range("X20").Select() | Out-Null
